$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new rows before the existing row 324 (2019-11-29), shifting
# all subsequent rows down by 9 (old last row 396 -> new last row 405).
$ws.Rows("324:332").Insert()

# New historical rows: timestamp, date, id, name, open, high, low, close, vol
$newRows = @(
    @(1574035200, "2019-11-18", 0.635, 0.645, 0.625, 0.64, 447000),
    @(1574121600, "2019-11-19", 0.64,  0.64,  0.625, 0.64, 245600),
    @(1574208000, "2019-11-20", 0.63,  0.63,  0.615, 0.615, 728000),
    @(1574294400, "2019-11-21", 0.62,  0.64,  0.62,  0.63, 606500),
    @(1574380800, "2019-11-22", 0.63,  0.63,  0.62,  0.63, 423400),
    @(1574640000, "2019-11-25", 0.635, 0.635, 0.615, 0.62, 278600),
    @(1574726400, "2019-11-26", 0.62,  0.625, 0.61,  0.625, 221800),
    @(1574812800, "2019-11-27", 0.615, 0.625, 0.605, 0.62, 369700),
    @(1574899200, "2019-11-28", 0.575, 0.58,  0.555, 0.5600000000000001, 2956300)
)

$startRow = 324
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Range("A$r").Value = $row[0]
    $ws.Range("E$r").Value = $row[2]
    $ws.Range("F$r").Value = $row[3]
    $ws.Range("G$r").Value = $row[4]
    $ws.Range("H$r").Value = $row[5]
    $ws.Range("I$r").Value = $row[6]
}

# Column B holds plain text dates ("YYYY-MM-DD") -- force text so Excel
# does not convert them into date serial numbers.
$ws.Range("B324:B332").NumberFormat = "@"
$ws.Range("B324").Value = "2019-11-18"
$ws.Range("B325").Value = "2019-11-19"
$ws.Range("B326").Value = "2019-11-20"
$ws.Range("B327").Value = "2019-11-21"
$ws.Range("B328").Value = "2019-11-22"
$ws.Range("B329").Value = "2019-11-25"
$ws.Range("B330").Value = "2019-11-26"
$ws.Range("B331").Value = "2019-11-27"
$ws.Range("B332").Value = "2019-11-28"

# Column C ("0201") would otherwise be parsed as the number 201 because it
# looks numeric; force text formatting so the leading zero is preserved.
$ws.Range("C324:C332").NumberFormat = "@"
$ws.Range("C324:C332").Value = "0201"

# Column D is a plain non-numeric string, no special handling required.
$ws.Range("D324:D332").Value = "NOVA"
